$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.801.62'
$ws.Range("E2").Value = '  -0.64%  '
$ws.Range("D3").Value = '1.932.97'
$ws.Range("E3").Value = '  -1.08%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.82'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4893'
$ws.Range("E7").Value = '  -0.32%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2955'
$ws.Range("E8").Value = '  -0.20%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06860'
$ws.Range("E9").Value = '  +0.45%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.24'
$ws.Range("E10").Value = '  +0.42%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '105.59'
$ws.Range("E11").Value = '  -1.47%  '
$ws.Range("D12").Value = '1.935.24'
$ws.Range("E12").Value = '  -0.94%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07781'
$ws.Range("E13").Value = '  -0.16%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.327'
$ws.Range("E14").Value = '  -2.23%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.7000'
$ws.Range("E15").Value = '  -0.69%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '272.94'
$ws.Range("E16").Value = '  -3.23%  '
$ws.Range("D17").Value = '30.816.69'
$ws.Range("E17").Value = '  -0.72%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007694'
$ws.Range("E18").Value = '  +0.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.614'
$ws.Range("E19").Value = '  +2.40%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.03'
$ws.Range("E20").Value = '  -1.34%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.508'
$ws.Range("E23").Value = '  +0.58%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.823'
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '164.63'
$ws.Range("E25").Value = '  -3.04%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '19.54'
$ws.Range("E26").Value = '  -2.23%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.160'
$ws.Range("E27").Value = '  -1.86%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1034'
$ws.Range("E28").Value = '  -2.18%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.387'
$ws.Range("E29").Value = '  -1.92%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.559'
$ws.Range("E30").Value = '  -1.11%  '
$ws.Range("E31").Value = '  -1.96%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.382'
$ws.Range("E32").Value = '  -1.42%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04878'
$ws.Range("E33").Value = '  -1.43%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7563'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.145'
$ws.Range("E35").Value = '  -2.07%  '
$ws.Range("E36").Value = '  -0.02%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.709'
$ws.Range("E37").Value = '  -0.84%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01997'
$ws.Range("E38").Value = '  -0.66%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '79.59'
$ws.Range("E39").Value = '  +7.76%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.655'
$ws.Range("E40").Value = '  -1.83%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.473'
$ws.Range("E41").Value = '  -1.82%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.064'
$ws.Range("E42").Value = '  -3.31%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8876'
$ws.Range("E43").Value = '  +0.34%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4438'
$ws.Range("E44").Value = '  -0.88%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '108.17'
$ws.Range("E45").Value = '  -1.13%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.871'
$ws.Range("E46").Value = '  -3.10%  '
$ws.Range("E47").Value = '  -0.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '982.61'
$ws.Range("E48").Value = '  -0.66%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1243'
$ws.Range("E49").Value = '  -1.73%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.15'
$ws.Range("E50").Value = '  +1.34%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.271'
$ws.Range("E51").Value = '  -0.45%  '
